$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header row (row 1): To do / Description / Status, bold
$ws.Range("B1").Value = "To do"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Status"
$ws.Range("B1:D1").Font.Bold = $true

# Add sequential numbering in column A for the existing 13 task rows (rows 2-14)
for ($i = 2; $i -le 14; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Move the active selection to B3 (matches the authored workbook state)
$null = $ws.Range("B3").Select()

# Set the page to portrait orientation
$ws.PageSetup.Orientation = 1

Write-Host "done"
